$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.197.19"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "2.033.30"
$ws.Range("E3").Value = "  +3.80%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'247.88"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").Value = "'0.631"
$ws.Range("E6").Value = "  +2.46%  "
$ws.Range("D7").Value = "'60.49"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.398"
$ws.Range("E9").Value = "  +6.50%  "
$ws.Range("D10").Value = "'0.0810"
$ws.Range("E10").Value = "  +3.06%  "
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").Value = "'15.29"
$ws.Range("E12").Value = "  +7.36%  "
$ws.Range("D13").Value = "'0.864"
$ws.Range("E13").Value = "  +4.19%  "
$ws.Range("D14").Value = "2.334.88"
$ws.Range("E14").Value = "  +3.89%  "
$ws.Range("D15").Value = "'22.42"
$ws.Range("E15").Value = "  +2.81%  "
$ws.Range("D16").Value = "'5.52"
$ws.Range("E16").Value = "  +5.07%  "
$ws.Range("D17").Value = "2.029.52"
$ws.Range("E17").Value = "  +3.42%  "
$ws.Range("D18").Value = "37.120.38"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").Value = "'70.89"
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").Value = "0.0₃0869"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").Value = "'5.27"
$ws.Range("E21").Value = "  +4.34%  "
$ws.Range("D22").Value = "'231.38"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'2.51"
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E26").Value = "  +3.62%  "
$ws.Range("D27").Value = "'163.71"
$ws.Range("E27").Value = "  +1.73%  "
$ws.Range("E28").Value = "  -2.78%  "
$ws.Range("D29").Value = "'19.86"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").Value = "'1.40"
$ws.Range("E30").Value = "  +6.76%  "
$ws.Range("E31").Value = "  +2.90%  "
$ws.Range("D32").Value = "'4.86"
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("D33").Value = "'0.0666"
$ws.Range("E33").Value = "  +8.66%  "
$ws.Range("D34").Value = "'4.55"
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("D35").Value = "'2.50"
$ws.Range("E35").Value = "  +10.43%  "
$ws.Range("D36").Value = "'3.49"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  +2.32%  "
$ws.Range("D39").Value = "'5.42"
$ws.Range("D40").Value = "'0.0981"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("E42").Value = "  +2.86%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'17.09"
$ws.Range("E43").Value = "  +7.59%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0215"
$ws.Range("E44").Value = "  +2.64%  "
$ws.Range("D45").Value = "'92.88"
$ws.Range("E45").Value = "  +4.78%  "
$ws.Range("E46").Value = "  +4.18%  "
$ws.Range("D47").Value = "1.389.78"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("E48").Value = "  +6.43%  "
$ws.Range("D49").Value = "'2.17"
$ws.Range("E49").Value = "  +19.41%  "
$ws.Range("D50").Value = "'2.86"
$ws.Range("E50").Value = "  +1.32%  "
$ws.Range("D51").Value = "'46.61"
$ws.Range("E51").Value = "  +3.10%  "
